# Update "想去人数" (want-to-go count) values in column F for sheets
# "展览" and "全部类型", which hold duplicate data tables.
$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 1765
    15 = 57
    25 = 298
    30 = 588
    41 = 559
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
